$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Append the new data row (row 6) to "example_another_table".
$ws2.Range("A6").Value = 5
$ws2.Range("B6").Value = 'BREAK ALL Ñ ''\\\\ \\ // '' "''"''"'
$ws2.Range("C6").Value = 4

# Update the saved selection/view state on sheet2.
[void]$ws2.Range("E15").Select()
